{"js": "// Apply the three content edits described by the diff:\n// 1. Remove the date \"2024.07.31\" text from the first (header) paragraph,\n//    leaving only the leading spaces run.\n// 2. Change \"130\" -> \"150\" in \"Material: FR-4 TG\u2265130.\" line.\n// 3. Change \"Matte black\" -> \"Black\" in the solder-mask color line.\n//\n// (The remaining hunks in the source diff are pure run-splitting/merging\n// artifacts that leave the visible text unchanged, so they are not\n// re-created here.)\n\nconst body = context.document.body;\n\n// 1) Remove the date text \"2024.07.31\"\nconst dateResults = body.search(\"2024.07.31\", { matchCase: true });\ndateResults.load(\"text\");\nawait context.sync();\nfor (let i = 0; i < dateResults.items.length; i++) {\n  dateResults.items[i].delete();\n}\nawait context.sync();\n\n// 2) \"130\" -> \"150\" in the material spec line\nconst matResults = body.search(\"TG\u2265130.\", { matchCase: true });\nmatResults.load(\"text\");\nawait context.sync();\nfor (let i = 0; i < matResults.items.length; i++) {\n  matResults.items[i].insertText(\"TG\u2265150.\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 3) \"Matte black\" -> \"Black\"\nconst colorResults = body.search(\"Matte black\", { matchCase: true });\ncolorResults.load(\"text\");\nawait context.sync();\nfor (let i = 0; i < colorResults.items.length; i++) {\n  colorResults.items[i].insertText(\"Black\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Apply the three content edits described by the diff:\n# 1. Remove the date \"2024.07.31\" text from the first (header) paragraph,\n#    leaving only the leading spaces run.\n# 2. Change \"130\" -> \"150\" in \"Material: FR-4 TG>=130.\" line.\n# 3. Change \"Matte black\" -> \"Black\" in the solder-mask color line.\n#\n# (The remaining hunks in the source diff are pure run-splitting/merging\n# artifacts that leave the visible text unchanged, so they are not\n# re-created here.)\n\n$d = $word.ActiveDocument\n\n# 1) Remove the date text \"2024.07.31\"\n$rng1 = $d.Content\n$rng1.Find.Execute(\"2024.07.31\", $false, $false, $false, $false, $false, $true, 1, $false, \"\", 0) | Out-Null\nif ($rng1.Find.Found) {\n    $rng1.Delete()\n}\n\n# 2) \"130\" -> \"150\" in the material spec line\n$rng2 = $d.Content\n$rng2.Find.Execute(\"130\", $false, $false, $false, $false, $false, $true, 1, $false, \"150\", 1) | Out-Null\n\n# 3) \"Matte black\" -> \"Black\"\n$rng3 = $d.Content\n$rng3.Find.Execute(\"Matte black\", $false, $false, $false, $false, $false, $true, 1, $false, \"Black\", 1) | Out-Null\n"}
